$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 11: May-V1.3 / 50K
$ws.Range("C11").Value = "May-V1.3"
$ws.Range("D11").Value = "50K"

# Rename "apr" to "aprV1.2" in C10
$ws.Range("C10").Value = "aprV1.2"

# Update selection to match target state
$ws.Range("C11").Select()
